$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price / 1h-volume data.
# Rows 44/45 also swap which coin (dogwifhat vs Fetch.AI) occupies each rank.
# A leading apostrophe forces plain-number-looking prices to stay text cells,
# matching how this sheet stores every Price/Volume value as text.

# Row 2
$ws.Range("D2").Value = '68.071.49'
$ws.Range("E2").Value = '  +2.44%  '
# Row 3
$ws.Range("D3").Value = '3.606.27'
$ws.Range("E3").Value = '  +1.13%  '
# Row 4
$ws.Range("E4").Value = '  -0.18%  '
# Row 5
$ws.Range("D5").Value = '''205.05'
$ws.Range("E5").Value = '  +10.77%  '
# Row 6
$ws.Range("D6").Value = '''564.76'
$ws.Range("E6").Value = '  -4.05%  '
# Row 7
$ws.Range("D7").Value = '3.601.04'
$ws.Range("E7").Value = '  +1.16%  '
# Row 8
$ws.Range("E8").Value = '  +1.20%  '
# Row 9
$ws.Range("E9").Value = '  -0.04%  '
# Row 10
$ws.Range("E10").Value = '  +0.26%  '
# Row 11
$ws.Range("D11").Value = '''60.84'
$ws.Range("E11").Value = '  +14.44%  '
# Row 12
$ws.Range("E12").Value = '  +3.11%  '
# Row 13
$ws.Range("D13").Value = '''0.0000287'
$ws.Range("E13").Value = '  +10.07%  '
# Row 14
$ws.Range("D14").Value = '''10.03'
$ws.Range("E14").Value = '  +2.06%  '
# Row 15
$ws.Range("D15").Value = '4.188.88'
$ws.Range("E15").Value = '  +1.24%  '
# Row 16
$ws.Range("D16").Value = '3.604.88'
$ws.Range("E16").Value = '  +1.21%  '
# Row 17
$ws.Range("E17").Value = '  +0.76%  '
# Row 18
$ws.Range("D18").Value = '''18.91'
$ws.Range("E18").Value = '  +3.07%  '
# Row 19
$ws.Range("D19").Value = '67.869.81'
$ws.Range("E19").Value = '  +2.39%  '
# Row 20
$ws.Range("D20").Value = '''12.37'
$ws.Range("E20").Value = '  +1.00%  '
# Row 21
$ws.Range("D21").Value = '''1.08'
$ws.Range("E21").Value = '  +1.91%  '
# Row 22
$ws.Range("D22").Value = '''402.32'
$ws.Range("E22").Value = '  +1.51%  '
# Row 23
$ws.Range("D23").Value = '''13.16'
$ws.Range("E23").Value = '  +16.59%  '
# Row 24
$ws.Range("E24").Value = '  -4.90%  '
# Row 25
$ws.Range("D25").Value = '''85.37'
$ws.Range("E25").Value = '  -0.67%  '
# Row 26
$ws.Range("D26").Value = '''4.01'
$ws.Range("E26").Value = '  +13.35%  '
# Row 27
$ws.Range("D27").Value = '''2.93'
$ws.Range("E27").Value = '  +1.09%  '
# Row 28
$ws.Range("D28").Value = '''12.60'
$ws.Range("E28").Value = '  +1.16%  '
# Row 29
$ws.Range("D29").Value = '''6.12'
$ws.Range("E29").Value = '  +1.24%  '
# Row 30
$ws.Range("D30").Value = '''8.33'
$ws.Range("E30").Value = '  +17.35%  '
# Row 31
$ws.Range("D31").Value = '''9.40'
$ws.Range("E31").Value = '  +4.78%  '
# Row 32
$ws.Range("D32").Value = '''31.63'
$ws.Range("E32").Value = '  +1.57%  '
# Row 33
$ws.Range("D33").Value = '''674.35'
$ws.Range("E33").Value = '  +8.26%  '
# Row 34
$ws.Range("D34").Value = '''12.23'
$ws.Range("E34").Value = '  +0.25%  '
# Row 35
$ws.Range("E35").Value = '  +0.68%  '
# Row 36
$ws.Range("D36").Value = '''63.89'
$ws.Range("E36").Value = '  +1.19%  '
# Row 37
$ws.Range("D37").Value = '''42.45'
$ws.Range("E37").Value = '  +2.69%  '
# Row 38
$ws.Range("D38").Value = '''0.422'
$ws.Range("E38").Value = '  +5.93%  '
# Row 39
$ws.Range("E39").Value = '  +0.01%  '
# Row 40
$ws.Range("E40").Value = '  -0.18%  '
# Row 41
$ws.Range("D41").Value = '3.297.42'
$ws.Range("E41").Value = '  +9.52%  '
# Row 42
$ws.Range("D42").Value = '''3.18'
$ws.Range("E42").Value = '  +12.99%  '
# Row 43
$ws.Range("E43").Value = '  +3.86%  '
# Row 44
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '''2.77'
$ws.Range("E44").Value = '  +9.57%  '
# Row 45
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''3.04'
$ws.Range("E45").Value = '  +29.59%  '
# Row 46
$ws.Range("D46").Value = '''0.998'
$ws.Range("E46").Value = '  -0.01%  '
# Row 47
$ws.Range("E47").Value = '  +2.14%  '
# Row 48
$ws.Range("D48").Value = '''2.75'
$ws.Range("E48").Value = '  +11.55%  '
# Row 49
$ws.Range("D49").Value = '''8.86'
$ws.Range("E49").Value = '  +3.36%  '
# Row 50
$ws.Range("E50").Value = '  +0.63%  '
# Row 51
$ws.Range("D51").Value = '''3.12'
$ws.Range("E51").Value = '  +0.99%  '
